# WR_89894498_WeekEnding_080325.xlsx — weekly units report update
# - Report generated timestamp refreshed
# - Total billed amount / line item count refreshed (now $0 billed, 8 line items)
# - Two additional completed-unit line items inserted above the TOTAL row
# - All per-line "Pricing" values (and the TOTAL) zeroed out

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert two new detail rows right above the existing TOTAL row (22),
#    pushing TOTAL from row 22 down to row 24. Excel auto-extends the
#    A22:G22 merged "TOTAL" cell down to A24:G24 and copies row 21's
#    formatting onto the two freshly inserted rows.
# ---------------------------------------------------------------------------
$ws.Rows("22:23").Insert()

# The engine clones row 21's (shaded) style onto both new rows; row 22 needs
# to go back to the unshaded style used by the other "odd" data rows
# (16/18/20), so pull that formatting back down from row 20.
$ws.Range("A20:H20").Copy()
$ws.Range("A22:H22").PasteSpecial(-4122)
$ws.Range("A22:H22").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Refresh header / summary details
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 8

# ---------------------------------------------------------------------------
# 3) Row 16-18 keep the same unit rows, but pricing now reads 0
# ---------------------------------------------------------------------------
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0

# ---------------------------------------------------------------------------
# 4) Row 19 becomes "Point 01 / PLD-EYE-ARM"
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "Point 01"
$ws.Range("B19").Value = "PLD-EYE-ARM"
$ws.Range("C19").Value = "Inst"
$ws.Range("D19").Value = "PLD, Eyebolt Deadend, Arm"
$ws.Range("E19").Value = "EA"
$ws.Range("F19").Value = 2
$ws.Range("H19").Value = 0

# ---------------------------------------------------------------------------
# 5) Row 20 stays "POL-40-2" but its point number becomes "Point 03"
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "Point 03"
$ws.Range("B20").Value = "POL-40-2"
$ws.Range("C20").Value = "Inst"
$ws.Range("D20").Value = "Pole,40ft,Class 2"
$ws.Range("E20").Value = "EA"
$ws.Range("F20").Value = 1
$ws.Range("H20").Value = 0

# ---------------------------------------------------------------------------
# 6) Row 21 becomes "Point 01 / INS-15-D-S"
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Point 01"
$ws.Range("B21").Value = "INS-15-D-S"
$ws.Range("C21").Value = "Inst"
$ws.Range("D21").Value = "INS,15kV,Deadend,Polymer"
$ws.Range("E21").Value = "EA"
$ws.Range("F21").Value = 2
$ws.Range("H21").Value = 0

# ---------------------------------------------------------------------------
# 7) New row 22: "Point 01 / POL-40-2" (newly inserted line item)
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = "Point 01"
$ws.Range("B22").Value = "POL-40-2"
$ws.Range("C22").Value = "Inst"
$ws.Range("D22").Value = "Pole,40ft,Class 2"
$ws.Range("E22").Value = "EA"
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0

# ---------------------------------------------------------------------------
# 8) New row 23: "Point 05 / POL-45-2" (newly inserted line item)
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "Point 05"
$ws.Range("B23").Value = "POL-45-2"
$ws.Range("C23").Value = "Inst"
$ws.Range("D23").Value = "Pole,45ft,Class 2"
$ws.Range("E23").Value = "EA"
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = 0

# ---------------------------------------------------------------------------
# 9) TOTAL row, now at row 24
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "TOTAL"
$ws.Range("H24").Value = 0
